$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append a new "Commit 11" data block (rows 211-226) below the existing
# "Commit 10" block, mirroring its layout/formulas exactly, per the commit:
#   "Unrolled the inner loop reducing cycle count to 7422"
# ---------------------------------------------------------------------------

# Row 211: section label (plain, unstyled - like row 194 above it)
$ws.Range("A211").Value = "Commit 11"

# Row 212: "MARS Tool Output" ... "Calulations" header
$ws.Range("A212").Value = "MARS Tool Output"
$ws.Range("D212").Value = "Calulations"

# Row 214: "Instruction Statistics Tool" sub-header
$ws.Range("A214").Value = "Instruction Statistics Tool"

# Row 215: column headers
$ws.Range("A215").Value = "Instruction type"
$ws.Range("B215").Value = "Count"
$ws.Range("D215").Value = "Adjusted count"
$ws.Range("E215").Value = "CPI"
$ws.Range("F215").Value = "Total cycles"

# Row 216: ALU
$ws.Range("A216").Value = "ALU"
$ws.Range("B216").Value = 2009
$ws.Range("D216").Formula = "=B216"
$ws.Range("E216").Value = 1
$ws.Range("F216").Formula = "=D216*E216"

# Row 217: Jump
$ws.Range("A217").Value = "Jump"
$ws.Range("B217").Value = 6
$ws.Range("D217").Formula = "=B217"
$ws.Range("E217").Value = 1
$ws.Range("F217").Formula = "=D217*E217"

# Row 218: Branch
$ws.Range("A218").Value = "Branch"
$ws.Range("B218").Value = 357
$ws.Range("D218").Formula = "=B218"
$ws.Range("E218").Value = 2
$ws.Range("F218").Formula = "=D218*E218"

# Row 219: Memory (no adjusted-count/CPI/cycles cells)
$ws.Range("A219").Value = "Memory"
$ws.Range("B219").Value = 474

# Row 220: Other
$ws.Range("A220").Value = "Other"
$ws.Range("B220").Value = 181
$ws.Range("D220").Formula = "=B220-(B224+B225-B219)"
$ws.Range("E220").Value = 5
$ws.Range("F220").Formula = "=D220*E220"

# Row 222: "Data Cache Simulation Tool" sub-header
$ws.Range("A222").Value = "Data Cache Simulation Tool"

# Row 223: column headers
$ws.Range("A223").Value = "Access"
$ws.Range("B223").Value = "Count"

# Row 224: Cache hit
$ws.Range("A224").Value = "Cache hit"
$ws.Range("B224").Value = 484
$ws.Range("D224").Formula = "=B224"
$ws.Range("E224").Value = 2
$ws.Range("F224").Formula = "=D224*E224"

# Row 225: Cache miss
$ws.Range("A225").Value = "Cache miss"
$ws.Range("B225").Value = 82
$ws.Range("D225").Formula = "=B225"
$ws.Range("E225").Value = 40
$ws.Range("F225").Formula = "=D225*E225"

# Row 226: grand total
$ws.Range("F226").Formula = "=SUM(F216:F225)"

# ---------------------------------------------------------------------------
# Formatting: copy the exact cell styles from the analogous cells in the
# "Commit 10" block above (rows 194-209) so style indices/semantics match.
# ---------------------------------------------------------------------------
$fmt = -4122  # xlPasteFormats

$ws.Range("A195:D195").Copy()
$ws.Range("A212:D212").PasteSpecial($fmt)

$ws.Range("A197").Copy()
$ws.Range("A214").PasteSpecial($fmt)

$ws.Range("A198:B198").Copy()
$ws.Range("A215:B215").PasteSpecial($fmt)
$ws.Range("D198:F198").Copy()
$ws.Range("D215:F215").PasteSpecial($fmt)

$ws.Range("A205").Copy()
$ws.Range("A222").PasteSpecial($fmt)

$ws.Range("A206:B206").Copy()
$ws.Range("A223:B223").PasteSpecial($fmt)

$ws.Range("F209").Copy()
$ws.Range("F226").PasteSpecial($fmt)

$ws.Application.CutCopyMode = $false

# Rows 212 and 226 use the bold 14pt header font, matching the taller
# row height (18.75) already used on their analogues (195 and 209).
$ws.Rows.Item(212).RowHeight = 18.75
$ws.Rows.Item(226).RowHeight = 18.75

$ws.Application.Calculate()

# ---------------------------------------------------------------------------
# Update the view so the new total row is visible, matching the saved state.
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 205
$win.ScrollColumn = 1
$ws.Range("I221").Select()
